$d = $word.ActiveDocument

# --- Change 1: revision-history table date "January 11, 2010" -> "1/11/2010" ---
$d.Content.Find.Execute("January 11, 2010", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1/11/2010", 2)

# --- Change 2: footer page-number field cached text "3" -> "1" ---
$d.Sections(2).Footers(1).Range.Find.Execute("3", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1", 2)
